$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New 4D box draw result for 2/7/2025 goes on top; every existing data row
# (2/7 .. 15/6) shifts down by one row. Row-level formatting (the ht=60
# custom row heights) stays put on its row number, so we only need to move
# the *values*, not the rows themselves - write the new table from the
# bottom up so we never clobber a value before it has been copied down.

$ws.Range("A9").Value = "15/6/2025 (Sun)"
$ws.Range("B9").Value = "2 9 3 2`n0 2 8 9`n8 5 2 5`n6 7 4 1"
$ws.Range("C9").Value = "✅ Direct: 12/3547 (0.34%)`n✅ iBet: 12/195 (6.15%)"

$ws.Range("A8").Value = "18/6/2025 (Wed)"
$ws.Range("B8").Value = "2 1 2 1`n3 2 4 7`n0 5 6 9`n5 3 3 8"
$ws.Range("C8").Value = "✅ Direct: 9/3416 (0.26%)`n✅ iBet: 9/188 (4.79%)"

$ws.Range("A7").Value = "21/6/2025 (Sat)"
$ws.Range("B7").Value = "4 6 1 8`n8 1 0 4`n1 5 7 7`n0 2 9 3"
$ws.Range("C7").Value = "✅ Direct: 13/3814 (0.34%)`n✅ iBet: 13/208 (6.25%)"

$ws.Range("A6").Value = "22/6/2025 (Sun)"
$ws.Range("B6").Value = "4 1 3 7`n6 2 5 4`n0 4 2 8`n9 5 6 3"
$ws.Range("C6").Value = "✅ Direct: 11/4144 (0.27%)`n✅ iBet: 11/222 (4.95%)"

$ws.Range("A5").Value = "25/6/2025 (Wed)"
$ws.Range("B5").Value = "2 0 5 7`n6 2 8 8`n9 1 3 0`n7 6 1 4"
$ws.Range("C5").Value = "✅ Direct: 12/4302 (0.28%)`n✅ iBet: 12/226 (5.31%)"

$ws.Range("A4").Value = "28/6/2025 (Sat)"
$ws.Range("B4").Value = "3 4 6 0`n4 9 3 6`n1 5 2 7`n0 0 4 8"
$ws.Range("C4").Value = "✅ Direct: 14/3980 (0.35%)`n✅ iBet: 14/215 (6.51%)"

$ws.Range("A3").Value = "29/6/2025 (Sun)"
$ws.Range("B3").Value = "6 5 2 6`n4 0 7 9`n0 4 6 5`n1 3 0 8"
$ws.Range("C3").Value = "✅ Direct: 14/3980 (0.35%)`n✅ iBet: 14/215 (6.51%)"

$ws.Range("A2").Value = "2/7/2025 (Wed)"
$ws.Range("B2").Value = "6 5 2 6`n4 0 7 9`n0 4 6 5`n1 3 0 8"
$ws.Range("C2").Value = "✅ Direct: 14/3980 (0.35%)`n✅ iBet: 14/215 (6.51%)"

# Row 2 never had a fixed custom height - writing multi-line wrapped text
# into it makes Excel auto-expand it, so put it back to the (non-custom)
# standard height.
$ws.Rows(2).AutoFit()

# Row 12 previously only had a formatted (wrap-text) blank cell in column B;
# give it the matching blank formatted cell in column C too.
$ws.Range("C12").WrapText = $true

# The table grew by one row overall, so the very last filler row (37, a
# lone formatted blank cell in column B) now needs a row-38 counterpart.
$ws.Range("B38").WrapText = $true
